$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6586
$ws.Range("E2").Value = 410
$ws.Range("F2").Value = 414
$ws.Range("G2").Value = 313
$ws.Range("H2").Value = 231
$ws.Range("I2").Value = 184
$ws.Range("J2").Value = 47
$ws.Range("K2").Value = 7231
$ws.Range("L2").Value = 4485
$ws.Range("M2").Value = 2746
$ws.Range("N2").Value = 2270
$ws.Range("O2").Value = 476
$ws.Range("P2").Value = 249
$ws.Range("Q2").Value = 1010
$ws.Range("R2").Value = -1393
$ws.Range("S2").Value = 381
$ws.Range("T2").Value = 1279
$ws.Range("U2").Value = -268
$ws.Range("V2").Value = 3275
$ws.Range("W2").Value = 6.22
$ws.Range("X2").Value = 3.5
$ws.Range("Y2").Value = 8.359999999999999
$ws.Range("Z2").Value = 3.3
$ws.Range("AA2").Value = 163.34
$ws.Range("AB2").Value = 836.79
$ws.Range("AC2").Value = 3701
$ws.Range("AD2").Value = 16.29
$ws.Range("AE2").Value = 10150
$ws.Range("AF2").Value = 1.19
$ws.Range("AH2").Value = 0.66
$ws.Range("AI2").Value = 9.75
$ws.Range("AJ2").Value = 24279820

# Row 3
$ws.Range("D3").Value = 5660
$ws.Range("E3").Value = 322
$ws.Range("F3").Value = 322
$ws.Range("G3").Value = 202
$ws.Range("H3").Value = 76
$ws.Range("I3").Value = 54
$ws.Range("J3").Value = 23
$ws.Range("K3").Value = 6783
$ws.Range("L3").Value = 4084
$ws.Range("M3").Value = 2700
$ws.Range("N3").Value = 2264
$ws.Range("O3").Value = 435
$ws.Range("P3").Value = 249
$ws.Range("Q3").Value = 188
$ws.Range("R3").Value = -367
$ws.Range("S3").Value = 68
$ws.Range("T3").Value = 315
$ws.Range("U3").Value = -127
$ws.Range("V3").Value = 3130
$ws.Range("W3").Value = 5.7
$ws.Range("X3").Value = 1.35
$ws.Range("Y3").Value = 2.37
$ws.Range("Z3").Value = 1.09
$ws.Range("AA3").Value = 151.28
$ws.Range("AB3").Value = 845.12
$ws.Range("AC3").Value = 1080
$ws.Range("AD3").Value = 36.39
$ws.Range("AE3").Value = 10124
$ws.Range("AF3").Value = 0.78
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 1.02
$ws.Range("AI3").Value = 33.42
$ws.Range("AJ3").Value = 24279820

# Row 4
$ws.Range("D4").Value = 6032
$ws.Range("E4").Value = 563
$ws.Range("F4").Value = 563
$ws.Range("G4").Value = 501
$ws.Range("H4").Value = 358
$ws.Range("I4").Value = 307
$ws.Range("J4").Value = 51
$ws.Range("K4").Value = 7108
$ws.Range("L4").Value = 4060
$ws.Range("M4").Value = 3048
$ws.Range("N4").Value = 2566
$ws.Range("O4").Value = 482
$ws.Range("P4").Value = 249
$ws.Range("Q4").Value = 466
$ws.Range("R4").Value = -257
$ws.Range("S4").Value = -172
$ws.Range("T4").Value = 196
$ws.Range("U4").Value = 270
$ws.Range("V4").Value = 2957
$ws.Range("W4").Value = 9.34
$ws.Range("X4").Value = 5.93
$ws.Range("Y4").Value = 12.69
$ws.Range("Z4").Value = 5.15
$ws.Range("AA4").Value = 133.2
$ws.Range("AB4").Value = 961.21
$ws.Range("AC4").Value = 6163
$ws.Range("AD4").Value = 6.34
$ws.Range("AE4").Value = 11474
$ws.Range("AF4").Value = 0.68
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 1.28
$ws.Range("AI4").Value = 7.31
$ws.Range("AJ4").Value = 24279820

# Row 5
$ws.Range("D5").Value = 6132
$ws.Range("E5").Value = 460
$ws.Range("F5").Value = 460
$ws.Range("G5").Value = 395
$ws.Range("H5").Value = 266
$ws.Range("I5").Value = 239
$ws.Range("J5").Value = 28
$ws.Range("K5").Value = 7564
$ws.Range("L5").Value = 4253
$ws.Range("M5").Value = 3311
$ws.Range("N5").Value = 2725
$ws.Range("O5").Value = 586
$ws.Range("P5").Value = 249
$ws.Range("Q5").Value = 98
$ws.Range("R5").Value = -159
$ws.Range("S5").Value = 190
$ws.Range("T5").Value = 142
$ws.Range("U5").Value = -44
$ws.Range("V5").Value = 3171
$ws.Range("W5").Value = 7.5
$ws.Range("X5").Value = 4.34
$ws.Range("Y5").Value = 9.029999999999999
$ws.Range("Z5").Value = 3.63
$ws.Range("AA5").Value = 128.45
$ws.Range("AB5").Value = 1045.75
$ws.Range("AC5").Value = 960
$ws.Range("AD5").Value = 5.98
$ws.Range("AE5").Value = 12186
$ws.Range("AF5").Value = 0.47
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 1.74
$ws.Range("AI5").Value = 9.390000000000001
$ws.Range("AJ5").Value = 24279820

# Row 6
$ws.Range("D6").Value = 5679
$ws.Range("E6").Value = 253
$ws.Range("F6").Value = 253
$ws.Range("G6").Value = 233
$ws.Range("H6").Value = 149
$ws.Range("I6").Value = 144
$ws.Range("K6").Value = 8528
$ws.Range("L6").Value = 4986
$ws.Range("M6").Value = 3542
$ws.Range("N6").Value = 2848
$ws.Range("P6").Value = 249
$ws.Range("Q6").Value = 219
$ws.Range("R6").Value = -395
$ws.Range("S6").Value = 319
$ws.Range("T6").Value = 378
$ws.Range("U6").Value = -159
$ws.Range("V6").Value = 3707
$ws.Range("W6").Value = 4.46
$ws.Range("X6").Value = 2.62
$ws.Range("Y6").Value = 5.18
$ws.Range("Z6").Value = 1.85
$ws.Range("AA6").Value = 140.76
$ws.Range("AB6").Value = 1093.86
$ws.Range("AC6").Value = 580
$ws.Range("AD6").Value = 7.1
$ws.Range("AE6").Value = 12735
$ws.Range("AF6").Value = 0.32
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 2.43
$ws.Range("AI6").Value = 15.53
$ws.Range("AJ6").Value = 24279820

# Rows 7-9: remove all numeric/data columns (D:AJ), keep A, B, C untouched
$ws.Range("D7:AJ9").ClearContents()
